$d = $word.ActiveDocument

# --- Change 1: cosmetic run-merge for the research question (text unchanged) ---
$d.Content.Find.Execute(
    "What is an appropriate current to heat up the nichrome wire in PDMS?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What is an appropriate current to heat up the nichrome wire in PDMS?", 2)

# --- Change 2: cosmetic run-merge for "Temperature of PDMS" (text unchanged) ---
$d.Content.Find.Execute(
    "Temperature of PDMS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Temperature of PDMS", 2)

# --- Change 3: "Triplicate" -> "Put in fridge to reset temperature." ---
$d.Content.Find.Execute(
    "Triplicate",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Put in fridge to reset temperature.", 2)

# --- Change 4: "Run steps 3-5 four more times for 750, 1000, 1250, 1500mA."
#     -> "Repeat steps 3 – 5 three more times" ---
$d.Content.Find.Execute(
    "Run steps 3-5 four more times for 750, 1000, 1250, 1500mA.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Repeat steps 3 " + [char]0x2013 + " 5 three more times", 2)

# --- Change 5: insert a brand new list item right after that paragraph ---
# (InsertParagraphAfter duplicates the "Repeat steps..." paragraph's pPr, so the
#  new paragraph automatically keeps the ListParagraph/numId=4 numbering.)
$findRange = $d.Content
$findRange.Find.Execute("Repeat steps 3 " + [char]0x2013 + " 5 three more times", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$repeatPara = $findRange.Paragraphs(1)
$repeatPara.Range.InsertParagraphAfter()
$repeatPara.Next().Range.Text = "Run steps 3-6 two more times for 1000 and 1500mA."

# --- Change 6/7: add a new, plain (non-list) empty paragraph after "Compare results." ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.ListFormat.RemoveNumbers()
$lastPara.Range.Style = "Normal"

Write-Output "edit complete"
